$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3903
$ws.Range("C3").Value = 3712
$ws.Range("C4").Value = 2836
$ws.Range("C5").Value = 1940
$ws.Range("C6").Value = 1746
$ws.Range("C7").Value = 836
$ws.Range("C8").Value = 577
$ws.Range("C9").Value = 549
$ws.Range("C10").Value = 507
$ws.Range("C11").Value = 494
